# exception_list.xlsx -- "added scheduler and minor tweaks due to database woes"
#
# * bandit: new exception id appended (row 11), cursor moved to F17
# * shark:  two new exception ids appended (rows 3-4), cursor moved to G21,
#           and shark becomes the active/selected sheet (was trust)
# * turn off concurrent calculation (database woes -> avoid concurrent recalcs)

$wb = $excel.ActiveWorkbook

$bandit = $wb.Worksheets.Item("bandit")
$shark  = $wb.Worksheets.Item("shark")

# --- bandit sheet: append the new exception id, then move the selection ---
$bandit.Activate()
$bandit.Range("A11").Value = 220989
$bandit.Range("F17").Select()

# --- shark sheet: append the two new exception ids, then move the selection ---
$shark.Activate()
$shark.Range("A3").Value = 219944
$shark.Range("A4").Value = 37542
$shark.Range("G21").Select()

# shark ends up the active / selected tab (was trust before this edit)
$shark.Activate()

# --- database woes: disable concurrent/multithreaded calculation ---
try {
    $excel.MultiThreadedCalculation.Enabled = $false
} catch {
}
